$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "福龙马"
$ws.Range("B2").Value = "三六零"
$ws.Range("C2").Value = "山子高科"
$ws.Range("A3").Value = "平潭发展"
$ws.Range("B3").Value = "赛微电子"
$ws.Range("C3").Value = "平潭发展"
$ws.Range("B4").Value = "平潭发展"
$ws.Range("C4").Value = "万向钱潮"
$ws.Range("A5").Value = "山子高科"
$ws.Range("A6").Value = "天际股份"
$ws.Range("B6").Value = "山子高科"
$ws.Range("A7").Value = "多氟多"
$ws.Range("B7").Value = "多氟多"
$ws.Range("C7").Value = "多氟多"
$ws.Range("A8").Value = "万向钱潮"
$ws.Range("B8").Value = "闻泰科技"
$ws.Range("C8").Value = "振德医疗"
$ws.Range("A9").Value = "闻泰科技"
$ws.Range("B9").Value = "大众公用"
$ws.Range("C9").Value = "天际股份"
$ws.Range("B10").Value = "天际股份"
$ws.Range("C10").Value = "时空科技"
$ws.Range("A11").Value = "海峡创新"
$ws.Range("B11").Value = "海峡创新"
$ws.Range("C11").Value = "和而泰"
$ws.Range("A12").Value = "荣科科技"
$ws.Range("B12").Value = "万向钱潮"
$ws.Range("C12").Value = "闻泰科技"
$ws.Range("A13").Value = "和而泰"
$ws.Range("B13").Value = "东方财富"
$ws.Range("C13").Value = "大众公用"
$ws.Range("A14").Value = "赛微电子"
$ws.Range("C14").Value = "三花智控"
$ws.Range("A15").Value = "众生药业"
$ws.Range("B15").Value = "贵州茅台"
$ws.Range("C15").Value = "海峡创新"
$ws.Range("A16").Value = "粤 传 媒"
$ws.Range("B16").Value = "粤 传 媒"
$ws.Range("C16").Value = "神州信息"
$ws.Range("A17").Value = "神州信息"
$ws.Range("B17").Value = "海马汽车"
$ws.Range("C17").Value = "工业富联"
$ws.Range("A18").Value = "海马汽车"
$ws.Range("B18").Value = "和而泰"
$ws.Range("C18").Value = "粤传媒"
$ws.Range("A19").Value = "金山办公"
$ws.Range("B19").Value = "神州信息"
$ws.Range("C19").Value = "海马汽车"
$ws.Range("A20").Value = "福昕软件"
$ws.Range("B20").Value = "众生药业"
$ws.Range("C20").Value = "合富中国"
$ws.Range("A21").Value = "工业富联"
$ws.Range("B21").Value = "荣科科技"
$ws.Range("C21").Value = "亚太药业"
